$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.138469934463501
$ws.Range("B1").Value = 4.555659770965576
$ws.Range("C1").Value = 3.410183191299438
$ws.Range("D1").Value = 0.8973651528358459
$ws.Range("E1").Value = 0.4719249606132507
